$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the Id value in the row that used to be last (typo fix: 57 -> 66) ---
$ws.Cells.Item(67,1).Value = 66

# --- Create rows 68-90, copying formatting from existing template rows ---
# Normal rows (BOOL/INT, has an Alarm/H column) are cloned from row 67.
# REAL 'SP' rows (no H column) are cloned from row 65.
# REAL 'PV' rows (no H column, Alarm=True) are cloned from row 66.
$ws.Range("A67:L67").Copy($ws.Range("A68:L68"))
$ws.Range("A67:L67").Copy($ws.Range("A69:L69"))
$ws.Range("A67:L67").Copy($ws.Range("A70:L70"))
$ws.Range("A67:L67").Copy($ws.Range("A71:L71"))
$ws.Range("A67:L67").Copy($ws.Range("A72:L72"))
$ws.Range("A67:L67").Copy($ws.Range("A73:L73"))
$ws.Range("A67:L67").Copy($ws.Range("A74:L74"))
$ws.Range("A67:L67").Copy($ws.Range("A75:L75"))
$ws.Range("A65:G65").Copy($ws.Range("A76:G76"))
$ws.Range("I65:L65").Copy($ws.Range("I76:L76"))
$ws.Range("A66:G66").Copy($ws.Range("A77:G77"))
$ws.Range("I66:L66").Copy($ws.Range("I77:L77"))
$ws.Range("A67:L67").Copy($ws.Range("A78:L78"))
$ws.Range("A67:L67").Copy($ws.Range("A79:L79"))
$ws.Range("A67:L67").Copy($ws.Range("A80:L80"))
$ws.Range("A67:L67").Copy($ws.Range("A81:L81"))
$ws.Range("A67:L67").Copy($ws.Range("A82:L82"))
$ws.Range("A67:L67").Copy($ws.Range("A83:L83"))
$ws.Range("A67:L67").Copy($ws.Range("A84:L84"))
$ws.Range("A67:L67").Copy($ws.Range("A85:L85"))
$ws.Range("A67:L67").Copy($ws.Range("A86:L86"))
$ws.Range("A65:G65").Copy($ws.Range("A87:G87"))
$ws.Range("I65:L65").Copy($ws.Range("I87:L87"))
$ws.Range("A66:G66").Copy($ws.Range("A88:G88"))
$ws.Range("I66:L66").Copy($ws.Range("I88:L88"))
$ws.Range("A67:L67").Copy($ws.Range("A89:L89"))
$ws.Range("A67:L67").Copy($ws.Range("A90:L90"))

# --- Column A: Id ---
$ws.Cells.Item(68,1).Value = 67
$ws.Cells.Item(69,1).Value = 68
$ws.Cells.Item(70,1).Value = 69
$ws.Cells.Item(71,1).Value = 70
$ws.Cells.Item(72,1).Value = 71
$ws.Cells.Item(73,1).Value = 72
$ws.Cells.Item(74,1).Value = 73
$ws.Cells.Item(75,1).Value = 74
$ws.Cells.Item(76,1).Value = 75
$ws.Cells.Item(77,1).Value = 76
$ws.Cells.Item(78,1).Value = 77
$ws.Cells.Item(79,1).Value = 78
$ws.Cells.Item(80,1).Value = 79
$ws.Cells.Item(81,1).Value = 80
$ws.Cells.Item(82,1).Value = 81
$ws.Cells.Item(83,1).Value = 82
$ws.Cells.Item(84,1).Value = 83
$ws.Cells.Item(85,1).Value = 84
$ws.Cells.Item(86,1).Value = 85
$ws.Cells.Item(87,1).Value = 86
$ws.Cells.Item(88,1).Value = 87
$ws.Cells.Item(89,1).Value = 88
$ws.Cells.Item(90,1).Value = 89

# --- Column B: Name ---
# (filled column-by-column across the whole new block, same order the original author used,
#  so new shared-string entries are appended in the same sequence as the source workbook)
$ws.Cells.Item(68,2).Value = 'P2_AUTO'
$ws.Cells.Item(69,2).Value = 'P2_MANUAL'
$ws.Cells.Item(70,2).Value = 'P2_START'
$ws.Cells.Item(71,2).Value = 'P2_STOP'
$ws.Cells.Item(72,2).Value = 'P2_BLOCKADE'
$ws.Cells.Item(73,2).Value = 'P2_RUN_H'
$ws.Cells.Item(74,2).Value = 'P2_RUN_M'
$ws.Cells.Item(75,2).Value = 'P2_RUN_S'
$ws.Cells.Item(76,2).Value = 'P2_SP'
$ws.Cells.Item(77,2).Value = 'P2_PV'
$ws.Cells.Item(78,2).Value = 'P2_MODE'

# --- Column C: Source ---
$ws.Cells.Item(68,3).Value = 'DB9.DBX18.2'
$ws.Cells.Item(69,3).Value = 'DB9.DBX18.3'
$ws.Cells.Item(70,3).Value = 'DB9.DBX18.0'
$ws.Cells.Item(71,3).Value = 'DB9.DBX18.1'
$ws.Cells.Item(72,3).Value = 'DB9.DBX18.4'
$ws.Cells.Item(73,3).Value = 'DB9.DBW20'
$ws.Cells.Item(74,3).Value = 'DB9.DBW22'
$ws.Cells.Item(75,3).Value = 'DB9.DBW24'
$ws.Cells.Item(76,3).Value = 'DB9.DBD30'
$ws.Cells.Item(77,3).Value = 'DB9.DBD34'
$ws.Cells.Item(78,3).Value = 'DB9.DBX26.1'

# --- Column B: Name (second block) ---
$ws.Cells.Item(79,2).Value = 'P3_AUTO'
$ws.Cells.Item(80,2).Value = 'P3_MANUAL'
$ws.Cells.Item(81,2).Value = 'P3_START'
$ws.Cells.Item(82,2).Value = 'P2_STOP'
$ws.Cells.Item(83,2).Value = 'P3_BLOCKADE'
$ws.Cells.Item(84,2).Value = 'P3_RUN_H'
$ws.Cells.Item(85,2).Value = 'P3_RUN_M'
$ws.Cells.Item(86,2).Value = 'P3_RUN_S'
$ws.Cells.Item(87,2).Value = 'P3_SP'
$ws.Cells.Item(88,2).Value = 'P3_PV'
$ws.Cells.Item(89,2).Value = 'P3_MODE'

# --- Column C: Source (second block) ---
$ws.Cells.Item(79,3).Value = 'DB15.DBX18.2'
$ws.Cells.Item(80,3).Value = 'DB15.DBX18.3'
$ws.Cells.Item(81,3).Value = 'DB15.DBX18.0'
$ws.Cells.Item(82,3).Value = 'DB15.DBX18.1'
$ws.Cells.Item(83,3).Value = 'DB15.DBX18.4'
$ws.Cells.Item(84,3).Value = 'DB15.DBW20'
$ws.Cells.Item(85,3).Value = 'DB15.DBW22'
$ws.Cells.Item(86,3).Value = 'DB15.DBW24'
$ws.Cells.Item(87,3).Value = 'DB15.DBD30'
$ws.Cells.Item(88,3).Value = 'DB15.DBD34'
$ws.Cells.Item(89,3).Value = 'DB15.DBX26.1'

# --- Column B & C for the final TestInt row ---
$ws.Cells.Item(90,2).Value = 'TestInt'
$ws.Cells.Item(90,3).Value = 'DB2.DBW4'

# --- Columns D (Type), E (Comment), F (AlarmText) ---
# These only reference shared strings that already exist in the workbook, so fill order is not significant.
$ws.Cells.Item(68,4).Value = 'BOOL'
$ws.Cells.Item(68,5).Value = 'Send'
$ws.Cells.Item(68,6).Value = 'Auto mode signal from Scada'
$ws.Cells.Item(69,4).Value = 'BOOL'
$ws.Cells.Item(69,5).Value = 'Send'
$ws.Cells.Item(69,6).Value = 'Manual mode signal from Scada'
$ws.Cells.Item(70,4).Value = 'BOOL'
$ws.Cells.Item(70,5).Value = 'Send'
$ws.Cells.Item(70,6).Value = 'Start signal from Scada'
$ws.Cells.Item(71,4).Value = 'BOOL'
$ws.Cells.Item(71,5).Value = 'Send'
$ws.Cells.Item(71,6).Value = 'Stop signal from Scada'
$ws.Cells.Item(72,4).Value = 'BOOL'
$ws.Cells.Item(72,5).Value = 'None'
$ws.Cells.Item(72,6).Value = '0 - blockade unactive, 1 - blockade active'
$ws.Cells.Item(73,4).Value = 'INT'
$ws.Cells.Item(73,5).Value = 'None'
$ws.Cells.Item(73,6).Value = 'Running time - hours'
$ws.Cells.Item(74,4).Value = 'INT'
$ws.Cells.Item(74,5).Value = 'None'
$ws.Cells.Item(74,6).Value = 'Running time - minutes'
$ws.Cells.Item(75,4).Value = 'INT'
$ws.Cells.Item(75,5).Value = 'None'
$ws.Cells.Item(75,6).Value = 'Running time - seconds'
$ws.Cells.Item(76,4).Value = 'REAL'
$ws.Cells.Item(76,5).Value = 'Send'
$ws.Cells.Item(76,6).Value = 'Setpoint value [%]'
$ws.Cells.Item(77,4).Value = 'REAL'
$ws.Cells.Item(77,5).Value = 'None'
$ws.Cells.Item(77,6).Value = 'Acctual speed [%]'
$ws.Cells.Item(78,4).Value = 'BOOL'
$ws.Cells.Item(78,5).Value = 'None'
$ws.Cells.Item(78,6).Value = '0 - auto, 1 - manual (Scada)'
$ws.Cells.Item(79,4).Value = 'BOOL'
$ws.Cells.Item(79,5).Value = 'Send'
$ws.Cells.Item(79,6).Value = 'Auto mode signal from Scada'
$ws.Cells.Item(80,4).Value = 'BOOL'
$ws.Cells.Item(80,5).Value = 'Send'
$ws.Cells.Item(80,6).Value = 'Manual mode signal from Scada'
$ws.Cells.Item(81,4).Value = 'BOOL'
$ws.Cells.Item(81,5).Value = 'Send'
$ws.Cells.Item(81,6).Value = 'Start signal from Scada'
$ws.Cells.Item(82,4).Value = 'BOOL'
$ws.Cells.Item(82,5).Value = 'Send'
$ws.Cells.Item(82,6).Value = 'Stop signal from Scada'
$ws.Cells.Item(83,4).Value = 'BOOL'
$ws.Cells.Item(83,5).Value = 'None'
$ws.Cells.Item(83,6).Value = '0 - blockade unactive, 1 - blockade active'
$ws.Cells.Item(84,4).Value = 'INT'
$ws.Cells.Item(84,5).Value = 'None'
$ws.Cells.Item(84,6).Value = 'Running time - hours'
$ws.Cells.Item(85,4).Value = 'INT'
$ws.Cells.Item(85,5).Value = 'None'
$ws.Cells.Item(85,6).Value = 'Running time - minutes'
$ws.Cells.Item(86,4).Value = 'INT'
$ws.Cells.Item(86,5).Value = 'None'
$ws.Cells.Item(86,6).Value = 'Running time - seconds'
$ws.Cells.Item(87,4).Value = 'REAL'
$ws.Cells.Item(87,5).Value = 'Send'
$ws.Cells.Item(87,6).Value = 'Setpoint value [%]'
$ws.Cells.Item(88,4).Value = 'REAL'
$ws.Cells.Item(88,5).Value = 'None'
$ws.Cells.Item(88,6).Value = 'Acctual speed [%]'
$ws.Cells.Item(89,4).Value = 'BOOL'
$ws.Cells.Item(89,5).Value = 'None'
$ws.Cells.Item(89,6).Value = '0 - auto, 1 - manual (Scada)'
$ws.Cells.Item(90,4).Value = 'INT'
$ws.Cells.Item(90,5).Value = 'None'
$ws.Cells.Item(90,6).Value = 'Test'

# Columns G (Alarm), J (AlarmLimitMin), K (AlarmLimitMax) and L (Historian) are left untouched:
# they were already copied correctly (False/0/1/False, or True for the two PV rows) from the template rows.

# --- Update the worksheet view to match where the user finished editing ---
$ws.Range("A88:A90").Select()
try {
  $excel.ActiveWindow.ScrollRow = 76
  $excel.ActiveWindow.ScrollColumn = 1
} catch {}

Write-Host "Rows 68-90 added successfully."
